$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 732.4541049094046
$ws.Range("C2").Value = 502.2941070197915
$ws.Range("D2").Value = 2480.16338965301
